# OrchestrationDesign.pptx edit:
#   * slide 2 title: "SAM Design - Orchestration" -> "SAM Design <U+2013> Orchestration-Old"
#   * slide 2: drop the "bessController/dockerController/Ryu" note text box (id=3)
#   * drop slides 3-10, leaving only the title slide and the "SAM Design" slide

$p = $ppt.ActivePresentation

# --- Update the title on slide 2 -----------------------------------------
$s2 = $p.Slides.Item(2)

$titleShape = $null
if ($s2.Shapes.HasTitle) {
    $titleShape = $s2.Shapes.Title
} else {
    for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
        $cand = $s2.Shapes.Item($i)
        if ($cand.Type -eq 14) { # msoPlaceholder
            if ($cand.PlaceholderFormat.Type -eq 13) { # ppPlaceholderTitle
                $titleShape = $cand
                break
            }
        }
    }
}

$enDash = [char]0x2013
$titleShape.TextFrame.TextRange.Text = "SAM Design " + $enDash + " Orchestration-Old"

# --- Remove the "bessController/dockerController/Ryu" text box on slide 2
for ($i = $s2.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s2.Shapes.Item($i)
    $isTarget = $false
    if ($sh.Id -eq 3) {
        $isTarget = $true
    } elseif ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -like "*bessController*") {
        $isTarget = $true
    }
    if ($isTarget) {
        $sh.Delete()
        break
    }
}

# --- Drop slides 3 through 10 (only the title + "SAM Design" slides stay)
for ($i = $p.Slides.Count; $i -ge 3; $i--) {
    $p.Slides.Item($i).Delete()
}
